$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.968.99'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '3.153.56'
$ws.Range("E3").Value = '  +1.93%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.95'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '625.85'
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("E7").Value = '  +22.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.369'
$ws.Range("E8").Value = '  -2.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").Value = '3.150.70'
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("E11").Value = '  +10.35%  '
$ws.Range("E12").Value = '  +6.01%  '
$ws.Range("E13").Value = '  +6.24%  '
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '35.19'
$ws.Range("E15").Value = '  +6.35%  '
$ws.Range("D16").Value = '90.750.54'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '3.738.24'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").Value = '3.125.87'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("E19").Value = '  +8.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.57'
$ws.Range("E20").Value = '  +5.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '472.33'
$ws.Range("E21").Value = '  +8.36%  '
$ws.Range("E22").Value = '  -5.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.17'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.20'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.76'
$ws.Range("E25").Value = '  +13.55%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.91'
$ws.Range("E26").Value = '  +5.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.37'
$ws.Range("E27").Value = '  +4.62%  '
$ws.Range("D28").Value = '3.321.38'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.230'
$ws.Range("E30").Value = '  +61.40%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.33'
$ws.Range("E31").Value = '  +5.83%  '
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.163'
$ws.Range("E32").Value = '  -2.73%  '
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.23'
$ws.Range("E34").Value = '  +18.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '520.83'
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.95'
$ws.Range("E36").Value = '  +5.85%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.146'
$ws.Range("E37").Value = '  +5.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.05'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("E39").Value = '  -7.07%  '
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0930'
$ws.Range("E41").Value = '  +28.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.433'
$ws.Range("E42").Value = '  +17.12%  '
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("E45").Value = '  +5.39%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.729'
$ws.Range("E47").Value = '  +19.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.76'
$ws.Range("E48").Value = '  +12.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '150.03'
$ws.Range("E49").Value = '  +6.35%  '
$ws.Range("E50").Value = '  +10.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '45.40'
